# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for price cells that would otherwise be
# auto-coerced into numbers (matches source data which stores these as text).
foreach ($addr in @("D5", "D11", "D14", "D19", "D20", "D21", "D22", "D25", "D29", "D36", "D38", "D40", "D47", "D48", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.761.32'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '3.308.13'
$ws.Range("E3").Value = '  +2.02%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '607.76'
$ws.Range("E5").Value = '  +2.38%  '

$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '3.303.09'
$ws.Range("E8").Value = '  +2.15%  '

$ws.Range("E10").Value = '  +2.00%  '

$ws.Range("D11").Value = '5.56'
$ws.Range("E11").Value = '  +3.87%  '

$ws.Range("E12").Value = '  +1.23%  '

$ws.Range("E13").Value = '  +0.80%  '

$ws.Range("D14").Value = '34.98'
$ws.Range("E14").Value = '  +2.15%  '

$ws.Range("D15").Value = '3.853.52'
$ws.Range("E15").Value = '  +2.20%  '

$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("D17").Value = '3.307.50'
$ws.Range("E17").Value = '  +2.18%  '

$ws.Range("D18").Value = '63.836.68'
$ws.Range("E18").Value = '  +0.82%  '

$ws.Range("D19").Value = '6.88'
$ws.Range("E19").Value = '  +1.98%  '

$ws.Range("D20").Value = '480.37'
$ws.Range("E20").Value = '  +1.63%  '

$ws.Range("D21").Value = '14.04'
$ws.Range("E21").Value = '  -0.67%  '

$ws.Range("D22").Value = '0.742'
$ws.Range("E22").Value = '  +1.59%  '

$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("E24").Value = '  +5.78%  '

$ws.Range("D25").Value = '85.35'
$ws.Range("E25").Value = '  +1.94%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("E27").Value = '  +1.81%  '

$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("D29").Value = '7.25'
$ws.Range("E29").Value = '  -1.64%  '

$ws.Range("E30").Value = '  +1.29%  '

$ws.Range("E31").Value = '  +1.34%  '

$ws.Range("E32").Value = '  +5.20%  '

$ws.Range("E33").Value = '  -0.45%  '

$ws.Range("E34").Value = '  +0.22%  '

$ws.Range("E35").Value = '  +0.97%  '

$ws.Range("D36").Value = '6.07'
$ws.Range("E36").Value = '  +2.68%  '

$ws.Range("D37").Value = '0.0₃0753'
$ws.Range("E37").Value = '  +6.17%  '

$ws.Range("D38").Value = '52.40'
$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("D40").Value = '432.42'
$ws.Range("E40").Value = '  +2.35%  '

$ws.Range("D41").Value = '3.114.12'
$ws.Range("E41").Value = '  +4.92%  '

$ws.Range("E42").Value = '  +8.56%  '

$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("E46").Value = '  +3.08%  '

$ws.Range("D47").Value = '36.85'
$ws.Range("E47").Value = '  +9.64%  '

$ws.Range("D48").Value = '26.42'
$ws.Range("E48").Value = '  +2.15%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '125.82'
$ws.Range("E50").Value = '  +3.54%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '2.32'
$ws.Range("E51").Value = '  -0.27%  '
